$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three pairing-minute entries (each appears twice, symmetric
# across the diagonal of the pairing matrix) from 240 to 330 minutes.
$ws.Range("F14").Value = 330
$ws.Range("B18").Value = 330

$ws.Range("E16").Value = 330
$ws.Range("D17").Value = 330

$ws.Range("H19").Value = 330
$ws.Range("G20").Value = 330

$ws.Range("J24").Select()
